$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Move row 2's data down to row 6 (rows 3-6 shift up by one) ---
$a2 = $ws.Cells.Item(2,1).Value()
$b2 = $ws.Cells.Item(2,2).Value()
$a3 = $ws.Cells.Item(3,1).Value()
$b3 = $ws.Cells.Item(3,2).Value()
$a4 = $ws.Cells.Item(4,1).Value()
$b4 = $ws.Cells.Item(4,2).Value()
$a5 = $ws.Cells.Item(5,1).Value()
$b5 = $ws.Cells.Item(5,2).Value()
$a6 = $ws.Cells.Item(6,1).Value()
$b6 = $ws.Cells.Item(6,2).Value()

$ws.Cells.Item(2,1).Value = $a3
$ws.Cells.Item(2,2).Value = $b3
$ws.Cells.Item(3,1).Value = $a4
$ws.Cells.Item(3,2).Value = $b4
$ws.Cells.Item(4,1).Value = $a5
$ws.Cells.Item(4,2).Value = $b5
$ws.Cells.Item(5,1).Value = $a6
$ws.Cells.Item(5,2).Value = $b6
$ws.Cells.Item(6,1).Value = $a2
$ws.Cells.Item(6,2).Value = $b2

# --- Stash the current "Hyperlink" and bordered "Normal" cell formats before they get disturbed ---
$ws.Range("A2").Copy()
$ws.Range("Z1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A4").Copy()
$ws.Range("Z2").PasteSpecial(-4122)   # xlPasteFormats
$ws.Application.CutCopyMode = $false

# --- Rebuild hyperlinks to match the new row positions ---
$ws.Range("A1").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("A6"), "mailto:sasikala.ars@gmail.com", "", "", "")

$keepA3 = $ws.Range("A3").Value()
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:sasikala.ars@gmail.com", "", "", "sasikala.ars@gmail.com")
$ws.Range("A3").Value = $keepA3

$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:sasikala.ars@gmail.com", "", "", "")

# --- Restore the correct cell formats (Add() re-stamps its own hyperlink look) ---
$ws.Range("Z1").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("A6").PasteSpecial(-4122)

$ws.Range("Z2").Copy()
$ws.Range("A3").PasteSpecial(-4122)

$ws.Application.CutCopyMode = $false
$ws.Range("Z1").Clear()
$ws.Range("Z2").Clear()

# --- Update the selected cell ---
$ws.Range("A3").Select()
